# 231030 - Pactos.xlsx: refresh "pacto" (B) / "partido" (C) mapping for the 2016 and 2021
# election blocks and add the missing incumbent ("XS"/"K" 2021) rows while dropping the
# now-obsolete tail rows (old 178-191), shrinking the used range to A1:C177.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the rows that no longer exist once the 2021 block is renumbered/merged.
$ws.Range("A178:A191").EntireRow.Delete()

# Rewrite rows 106-177 (year / pacto / partido) with their corrected values.
$ws.Cells.Item(106, 1).Value = 2016
$ws.Cells.Item(106, 2).Value = 'A'
$ws.Cells.Item(106, 3).Value = 'PARTIDO REGIONALISTA DE MAGALLANES'
$ws.Cells.Item(107, 1).Value = 2016
$ws.Cells.Item(107, 2).Value = 'C'
$ws.Cells.Item(107, 3).Value = 'PODER'
$ws.Cells.Item(108, 1).Value = 2016
$ws.Cells.Item(108, 2).Value = 'C'
$ws.Cells.Item(108, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(109, 1).Value = 2016
$ws.Cells.Item(109, 2).Value = 'C'
$ws.Cells.Item(109, 3).Value = 'PARTIDO ECOLOGISTA VERDE'
$ws.Cells.Item(110, 1).Value = 2016
$ws.Cells.Item(110, 2).Value = 'D'
$ws.Cells.Item(110, 3).Value = 'SOMOS AYSEN'
$ws.Cells.Item(111, 1).Value = 2016
$ws.Cells.Item(111, 2).Value = 'E'
$ws.Cells.Item(111, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(112, 1).Value = 2016
$ws.Cells.Item(112, 2).Value = 'E'
$ws.Cells.Item(112, 3).Value = 'PARTIDO POR LA DEMOCRACIA'
$ws.Cells.Item(113, 1).Value = 2016
$ws.Cells.Item(113, 2).Value = 'E'
$ws.Cells.Item(113, 3).Value = 'PARTIDO SOCIALISTA DE CHILE'
$ws.Cells.Item(114, 1).Value = 2016
$ws.Cells.Item(114, 2).Value = 'E'
$ws.Cells.Item(114, 3).Value = 'PARTIDO DEMOCRATA CRISTIANO'
$ws.Cells.Item(115, 1).Value = 2016
$ws.Cells.Item(115, 2).Value = 'E'
$ws.Cells.Item(115, 3).Value = 'PARTIDO RADICAL SOCIALDEMOCRATA'
$ws.Cells.Item(116, 1).Value = 2016
$ws.Cells.Item(116, 2).Value = 'E'
$ws.Cells.Item(116, 3).Value = 'MAS REGION'
$ws.Cells.Item(117, 1).Value = 2016
$ws.Cells.Item(117, 2).Value = 'E'
$ws.Cells.Item(117, 3).Value = 'PARTIDO COMUNISTA DE CHILE'
$ws.Cells.Item(118, 1).Value = 2016
$ws.Cells.Item(118, 2).Value = 'F'
$ws.Cells.Item(118, 3).Value = 'UNION DEMOCRATA INDEPENDIENTE'
$ws.Cells.Item(119, 1).Value = 2016
$ws.Cells.Item(119, 2).Value = 'F'
$ws.Cells.Item(119, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(120, 1).Value = 2016
$ws.Cells.Item(120, 2).Value = 'F'
$ws.Cells.Item(120, 3).Value = 'RENOVACION NACIONAL'
$ws.Cells.Item(121, 1).Value = 2016
$ws.Cells.Item(121, 2).Value = 'F'
$ws.Cells.Item(121, 3).Value = 'PARTIDO REGIONALISTA INDEPENDIENTE'
$ws.Cells.Item(122, 1).Value = 2016
$ws.Cells.Item(122, 2).Value = 'I'
$ws.Cells.Item(122, 3).Value = 'AMPLITUD'
$ws.Cells.Item(123, 1).Value = 2016
$ws.Cells.Item(123, 2).Value = 'I'
$ws.Cells.Item(123, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(124, 1).Value = 2016
$ws.Cells.Item(124, 2).Value = 'K'
$ws.Cells.Item(124, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(125, 1).Value = 2016
$ws.Cells.Item(125, 2).Value = 'K'
$ws.Cells.Item(125, 3).Value = 'REVOLUCION DEMOCRATICA'
$ws.Cells.Item(126, 1).Value = 2016
$ws.Cells.Item(126, 2).Value = 'M'
$ws.Cells.Item(126, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(127, 1).Value = 2016
$ws.Cells.Item(127, 2).Value = 'M'
$ws.Cells.Item(127, 3).Value = 'PARTIDO IGUALDAD'
$ws.Cells.Item(128, 1).Value = 2016
$ws.Cells.Item(128, 2).Value = 'M'
$ws.Cells.Item(128, 3).Value = 'PARTIDO FRENTE POPULAR'
$ws.Cells.Item(129, 1).Value = 2016
$ws.Cells.Item(129, 2).Value = 'N'
$ws.Cells.Item(129, 3).Value = 'FUERZA REGIONAL NORTE VERDE'
$ws.Cells.Item(130, 1).Value = 2016
$ws.Cells.Item(130, 2).Value = 'O'
$ws.Cells.Item(130, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(131, 1).Value = 2016
$ws.Cells.Item(131, 2).Value = 'O'
$ws.Cells.Item(131, 3).Value = 'DEMOCRACIA REGIONAL PATAGONICA'
$ws.Cells.Item(132, 1).Value = 2016
$ws.Cells.Item(132, 2).Value = 'O'
$ws.Cells.Item(132, 3).Value = 'PARTIDO PROGRESISTA'
$ws.Cells.Item(133, 1).Value = 2016
$ws.Cells.Item(133, 2).Value = 'O'
$ws.Cells.Item(133, 3).Value = 'FRENTE REGIONAL Y POPULAR'
$ws.Cells.Item(134, 1).Value = 2016
$ws.Cells.Item(134, 2).Value = 'O'
$ws.Cells.Item(134, 3).Value = 'WALLMAPUWEN'
$ws.Cells.Item(135, 1).Value = 2016
$ws.Cells.Item(135, 2).Value = 'P'
$ws.Cells.Item(135, 3).Value = 'PARTIDO HUMANISTA'
$ws.Cells.Item(136, 1).Value = 2016
$ws.Cells.Item(136, 2).Value = 'P'
$ws.Cells.Item(136, 3).Value = 'PARTIDO LIBERAL DE CHILE'
$ws.Cells.Item(137, 1).Value = 2016
$ws.Cells.Item(137, 2).Value = 'P'
$ws.Cells.Item(137, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(138, 1).Value = 2016
$ws.Cells.Item(138, 2).Value = 'P'
$ws.Cells.Item(138, 3).Value = 'MOVIMIENTO INDEPENDIENTE REGIONALISTA AGRARIO Y SOCIAL'
$ws.Cells.Item(139, 1).Value = 2016
$ws.Cells.Item(139, 2).Value = 'Q'
$ws.Cells.Item(139, 3).Value = 'UNIDOS RESULTA EN DEMOCRACIA'
$ws.Cells.Item(140, 1).Value = 2016
$ws.Cells.Item(140, 2).Value = 'R'
$ws.Cells.Item(140, 3).Value = 'UNION PATRIOTICA'
$ws.Cells.Item(141, 1).Value = 2016
$ws.Cells.Item(141, 2).Value = 'R'
$ws.Cells.Item(141, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(142, 1).Value = 2016
$ws.Cells.Item(142, 2).Value = 'T'
$ws.Cells.Item(142, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(143, 1).Value = 2021
$ws.Cells.Item(143, 2).Value = 'CANDIDATURA INDEPENDIENTE'
$ws.Cells.Item(143, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(144, 1).Value = 2021
$ws.Cells.Item(144, 2).Value = 'K'
$ws.Cells.Item(144, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(145, 1).Value = 2021
$ws.Cells.Item(145, 2).Value = 'K'
$ws.Cells.Item(145, 3).Value = 'PARTIDO POR LA DEMOCRACIA'
$ws.Cells.Item(146, 1).Value = 2021
$ws.Cells.Item(146, 2).Value = 'K'
$ws.Cells.Item(146, 3).Value = 'PARTIDO RADICAL DE CHILE'
$ws.Cells.Item(147, 1).Value = 2021
$ws.Cells.Item(147, 2).Value = 'K'
$ws.Cells.Item(147, 3).Value = 'PARTIDO SOCIALISTA DE CHILE'
$ws.Cells.Item(148, 1).Value = 2021
$ws.Cells.Item(148, 2).Value = 'M'
$ws.Cells.Item(148, 3).Value = 'PARTIDO COMUNISTA DE CHILE'
$ws.Cells.Item(149, 1).Value = 2021
$ws.Cells.Item(149, 2).Value = 'M'
$ws.Cells.Item(149, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(150, 1).Value = 2021
$ws.Cells.Item(150, 2).Value = 'M'
$ws.Cells.Item(150, 3).Value = 'FEDERACION REGIONALISTA VERDE SOCIAL'
$ws.Cells.Item(151, 1).Value = 2021
$ws.Cells.Item(151, 2).Value = 'XE'
$ws.Cells.Item(151, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(152, 1).Value = 2021
$ws.Cells.Item(152, 2).Value = 'XE'
$ws.Cells.Item(152, 3).Value = 'PARTIDO NACIONAL CIUDADANO'
$ws.Cells.Item(153, 1).Value = 2021
$ws.Cells.Item(153, 2).Value = 'XO'
$ws.Cells.Item(153, 3).Value = 'PARTIDO REPUBLICANO DE CHILE'
$ws.Cells.Item(154, 1).Value = 2021
$ws.Cells.Item(154, 2).Value = 'XO'
$ws.Cells.Item(154, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(155, 1).Value = 2021
$ws.Cells.Item(155, 2).Value = 'XS'
$ws.Cells.Item(155, 3).Value = 'PARTIDO LIBERAL DE CHILE'
$ws.Cells.Item(156, 1).Value = 2021
$ws.Cells.Item(156, 2).Value = 'XS'
$ws.Cells.Item(156, 3).Value = 'CONVERGENCIA SOCIAL'
$ws.Cells.Item(157, 1).Value = 2021
$ws.Cells.Item(157, 2).Value = 'XS'
$ws.Cells.Item(157, 3).Value = 'REVOLUCION DEMOCRATICA'
$ws.Cells.Item(158, 1).Value = 2021
$ws.Cells.Item(158, 2).Value = 'XS'
$ws.Cells.Item(158, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(159, 1).Value = 2021
$ws.Cells.Item(159, 2).Value = 'XS'
$ws.Cells.Item(159, 3).Value = 'COMUNES'
$ws.Cells.Item(160, 1).Value = 2021
$ws.Cells.Item(160, 2).Value = 'XU'
$ws.Cells.Item(160, 3).Value = 'PARTIDO DEMOCRATA CRISTIANO'
$ws.Cells.Item(161, 1).Value = 2021
$ws.Cells.Item(161, 2).Value = 'XU'
$ws.Cells.Item(161, 3).Value = 'PARTIDO PROGRESISTA DE CHILE'
$ws.Cells.Item(162, 1).Value = 2021
$ws.Cells.Item(162, 2).Value = 'XU'
$ws.Cells.Item(162, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(163, 1).Value = 2021
$ws.Cells.Item(163, 2).Value = 'XU'
$ws.Cells.Item(163, 3).Value = 'CIUDADANOS'
$ws.Cells.Item(164, 1).Value = 2021
$ws.Cells.Item(164, 2).Value = 'XX'
$ws.Cells.Item(164, 3).Value = 'UNION DEMOCRATA INDEPENDIENTE'
$ws.Cells.Item(165, 1).Value = 2021
$ws.Cells.Item(165, 2).Value = 'XX'
$ws.Cells.Item(165, 3).Value = 'RENOVACION NACIONAL'
$ws.Cells.Item(166, 1).Value = 2021
$ws.Cells.Item(166, 2).Value = 'XX'
$ws.Cells.Item(166, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(167, 1).Value = 2021
$ws.Cells.Item(167, 2).Value = 'XX'
$ws.Cells.Item(167, 3).Value = 'EVOLUCION POLITICA'
$ws.Cells.Item(168, 1).Value = 2021
$ws.Cells.Item(168, 2).Value = 'XX'
$ws.Cells.Item(168, 3).Value = 'PARTIDO REGIONALISTA INDEPENDIENTE DEMOCRATA'
$ws.Cells.Item(169, 1).Value = 2021
$ws.Cells.Item(169, 2).Value = 'XY'
$ws.Cells.Item(169, 3).Value = 'IGUALDAD'
$ws.Cells.Item(170, 1).Value = 2021
$ws.Cells.Item(170, 2).Value = 'XY'
$ws.Cells.Item(170, 3).Value = 'PARTIDO HUMANISTA'
$ws.Cells.Item(171, 1).Value = 2021
$ws.Cells.Item(171, 2).Value = 'XY'
$ws.Cells.Item(171, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(172, 1).Value = 2021
$ws.Cells.Item(172, 2).Value = 'XZ'
$ws.Cells.Item(172, 3).Value = 'PARTIDO ECOLOGISTA VERDE'
$ws.Cells.Item(173, 1).Value = 2021
$ws.Cells.Item(173, 2).Value = 'XZ'
$ws.Cells.Item(173, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(174, 1).Value = 2021
$ws.Cells.Item(174, 2).Value = 'YC'
$ws.Cells.Item(174, 3).Value = 'PARTIDO CONSERVADOR CRISTIANO'
$ws.Cells.Item(175, 1).Value = 2021
$ws.Cells.Item(175, 2).Value = 'YC'
$ws.Cells.Item(175, 3).Value = 'INDEPENDIENTE'
$ws.Cells.Item(176, 1).Value = 2021
$ws.Cells.Item(176, 2).Value = 'YG'
$ws.Cells.Item(176, 3).Value = 'NUEVO TIEMPO'
$ws.Cells.Item(177, 1).Value = 2021
$ws.Cells.Item(177, 2).Value = 'ZB'
$ws.Cells.Item(177, 3).Value = 'UNION PATRIOTICA'
